$wb = $excel.ActiveWorkbook

# Remember what was selected/active before we start, so we can restore it
# (our edit only touches the "Edit Repayment Schedule" sheet).
$origSheet = $wb.ActiveSheet
$origSelection = $excel.Selection.Address()

$ws = $wb.Worksheets.Item("Edit Repayment Schedule")
$ws.Activate() | Out-Null

# Insert a new row above the current row 6 ("clickonsubmit"/"Submit"),
# shifting it and all rows below it down by one.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new "wait to page load" step.
$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# B6 should look like the other "wait" row (B3, green-filled amount cell).
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B6").Value = 2000

# Match the selection recorded in the edited workbook for this sheet.
$ws.Range("A6:B6").Select() | Out-Null

# Restore the originally active sheet/selection.
$origSheet.Activate() | Out-Null
$origSheet.Range($origSelection).Select() | Out-Null
